$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.8882575757575758
$ws.Range("I2").Value = 0.8882575757575758
$ws.Range("L2").Value = 24
$ws.Range("M2").Value = 48
$ws.Range("N2").Value = 0.6161616161616161
$ws.Range("T2").Value = 29
$ws.Range("U2").Value = 33
$ws.Range("V2").Value = 0.7627118644067796
$ws.Range("W2").Value = 0.7260869565217392
$ws.Range("Y2").Value = 0.7521739130434782
$ws.Range("Z2").Value = 0.7304347826086957
$ws.Range("AA2").Value = 0.7260869565217392
